$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.848.69'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '2.405.07'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.72'
$ws.Range('D6').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('D9').Value = '2.384.53'
$ws.Range('E9').Value = '  -2.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.06'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('E13').Value = '  -2.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.93'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').Value = '2.822.03'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '60.753.81'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '2.393.86'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.52'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.16'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.01'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.74%  '
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('E25').Value = '  -6.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.28'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '573.54'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -7.98%  '
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').Value = '0.0₃0908'
$ws.Range('E30').Value = '  -5.26%  '
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.34'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.84%  '
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  -8.46%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.60'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.62%  '
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.37'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '147.57'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.15'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.07'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.24%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.62'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.33'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.82%  '
$ws.Range('D46').Value = '0.0₆0282'
$ws.Range('E46').Value = '  +17.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '139.71'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.26'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.04%  '
